{"js": "const pairs = [\n  [\"2025-10-17 Friday\", \"2025-10-18 Saturday\"],\n  [\"55\u00d743=2365\", \"26\u00d720=520\"],\n  [\"41\u00d721=861\", \"16\u00d749=784\"],\n  [\"71\u00d733=2343\", \"29\u00d781=2349\"],\n  [\"66\u00d782=5412\", \"32\u00d773=2336\"],\n  [\"76\u00d764=4864\", \"49\u00d722=1078\"],\n  [\"92\u00d754=4968\", \"19\u00d768=1292\"],\n  [\"78\u00d764=4992\", \"20\u00d791=1820\"],\n  [\"54\u00d794=5076\", \"31\u00d711=341\"],\n  [\"27\u00d781=2187\", \"85\u00d725=2125\"],\n  [\"73\u00d716=1168\", \"18\u00d742=756\"],\n  [\"86\u00d729=2494\", \"42\u00d717=714\"],\n  [\"53\u00d786=4558\", \"40\u00d723=920\"],\n  [\"59\u00d789=5251\", \"65\u00d728=1820\"],\n  [\"74\u00d734=2516\", \"71\u00d743=3053\"],\n  [\"20\u00d738=760\", \"40\u00d782=3280\"],\n  [\"26\u00d796=2496\", \"61\u00d798=5978\"],\n  [\"33\u00d757=1881\", \"66\u00d717=1122\"],\n  [\"67\u00d711=737\", \"94\u00d789=8366\"],\n  [\"78\u00d732=2496\", \"60\u00d774=4440\"],\n  [\"32\u00d719=608\", \"48\u00d782=3936\"],\n  [\"21\u00d789=1869\", \"49\u00d742=2058\"],\n  [\"99\u00d741=4059\", \"26\u00d795=2470\"],\n  [\"34\u00d756=1904\", \"87\u00d771=6177\"],\n  [\"71\u00d758=4118\", \"61\u00d764=3904\"],\n  [\"46\u00d755=2530\", \"24\u00d745=1080\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    console.log(\"WARNING: text not found, skipped: \" + oldText);\n    continue;\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n$pairs = @(\n  ,@(\"2025-10-17 Friday\", \"2025-10-18 Saturday\")\n  ,@(\"55\u00d743=2365\", \"26\u00d720=520\")\n  ,@(\"41\u00d721=861\", \"16\u00d749=784\")\n  ,@(\"71\u00d733=2343\", \"29\u00d781=2349\")\n  ,@(\"66\u00d782=5412\", \"32\u00d773=2336\")\n  ,@(\"76\u00d764=4864\", \"49\u00d722=1078\")\n  ,@(\"92\u00d754=4968\", \"19\u00d768=1292\")\n  ,@(\"78\u00d764=4992\", \"20\u00d791=1820\")\n  ,@(\"54\u00d794=5076\", \"31\u00d711=341\")\n  ,@(\"27\u00d781=2187\", \"85\u00d725=2125\")\n  ,@(\"73\u00d716=1168\", \"18\u00d742=756\")\n  ,@(\"86\u00d729=2494\", \"42\u00d717=714\")\n  ,@(\"53\u00d786=4558\", \"40\u00d723=920\")\n  ,@(\"59\u00d789=5251\", \"65\u00d728=1820\")\n  ,@(\"74\u00d734=2516\", \"71\u00d743=3053\")\n  ,@(\"20\u00d738=760\", \"40\u00d782=3280\")\n  ,@(\"26\u00d796=2496\", \"61\u00d798=5978\")\n  ,@(\"33\u00d757=1881\", \"66\u00d717=1122\")\n  ,@(\"67\u00d711=737\", \"94\u00d789=8366\")\n  ,@(\"78\u00d732=2496\", \"60\u00d774=4440\")\n  ,@(\"32\u00d719=608\", \"48\u00d782=3936\")\n  ,@(\"21\u00d789=1869\", \"49\u00d742=2058\")\n  ,@(\"99\u00d741=4059\", \"26\u00d795=2470\")\n  ,@(\"34\u00d756=1904\", \"87\u00d771=6177\")\n  ,@(\"71\u00d758=4118\", \"61\u00d764=3904\")\n  ,@(\"46\u00d755=2530\", \"24\u00d745=1080\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  $find = $range.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $found = $find.Execute($oldText, $false, $true, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n  if (-not $found) {\n    Write-Host \"WARNING: text not found, skipped: $oldText\"\n  }\n}"}
